# Atualizado por script em 23-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the match data (columns F:V) between row 39 and row 40.
#    Columns A (Indice), B, C, D and E (data_partida) stay as-is in both
#    rows; only the match-specific data (home/away teams, goals, odds and
#    url) are exchanged between the two rows.
# ---------------------------------------------------------------------------
$ws.Range("F39:V39").Copy()
$ws.Range("F200:V200").PasteSpecial(-4104)   # xlPasteAll -> stash row 39 data

$ws.Range("F40:V40").Copy()
$ws.Range("F39:V39").PasteSpecial(-4104)     # row 40 data -> row 39

$ws.Range("F200:V200").Copy()
$ws.Range("F40:V40").PasteSpecial(-4104)     # stashed row 39 data -> row 40

$ws.Range("F200:V200").ClearContents()       # remove temporary stash

# ---------------------------------------------------------------------------
# 2) Append a new match as row 51 (Indice 50), copying the formatting from
#    row 50 and then filling in the new values.
# ---------------------------------------------------------------------------
$ws.Range("A50:V50").Copy()
$ws.Range("A51:V51").PasteSpecial(-4122)     # xlPasteFormats

$ws.Range("A51").Value = 50
$ws.Range("B51").Value = "united-arab-emirates"
$ws.Range("C51").Value = "uae-league"
$ws.Range("D51").Value = "2023-2024"
$ws.Range("E51").Value = 45253.6875
$ws.Range("F51").Value = "Al Bataeh"
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = "Al Sharjah"
$ws.Range("I51").Value = 3
$ws.Range("J51").Value = 4.95
$ws.Range("K51").Value = "18/11/2023 13:42"
$ws.Range("L51").Value = 5.43
$ws.Range("M51").Value = "23/11/2023 16:15"
$ws.Range("N51").Value = 4.42
$ws.Range("O51").Value = "18/11/2023 13:42"
$ws.Range("P51").Value = 4.56
$ws.Range("Q51").Value = "23/11/2023 16:15"
$ws.Range("R51").Value = 1.53
$ws.Range("S51").Value = "18/11/2023 13:42"
$ws.Range("T51").Value = 1.55
$ws.Range("U51").Value = "23/11/2023 16:15"
$ws.Range("V51").Value = "https://www.betexplorer.com/football/united-arab-emirates/uae-league/al-bataeh-al-sharjah/02DBKqeG/"

Write-Host "Edit complete"
